# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" summary text (cell A1) ---
$hoja1 = $wb.Worksheets.Item("Hoja1")
$hoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.08 = 7905.14 pesos`n✅ 7905.14 pesos = 2.09 = 938.6 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- tasas: update the rate cells N10/O10 and N12/O12 ---
$tasas = $wb.Worksheets.Item("tasas")
$tasas.Range("N10").Value = 480.7
$tasas.Range("O10").Value = 3800
$tasas.Range("N12").Value = 3790
$tasas.Range("O12").Value = 450
